$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.875.41'
$ws.Range("E2").Value = '  -1.44%  '
$ws.Range("D3").Value = '3.333.91'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.44'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.41%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.590'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.63%  '
$ws.Range("D9").Value = '3.330.85'
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.178'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.92%  '
$ws.Range("E11").Value = '  -1.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.70'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '659.67'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.86%  '
$ws.Range("D15").Value = '3.870.13'
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("E16").Value = '  -1.59%  '
$ws.Range("D17").Value = '67.985.56'
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("D19").Value = '3.332.77'
$ws.Range("E19").Value = '  -1.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.46'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.30%  '
$ws.Range("E22").Value = '  -2.34%  '
$ws.Range("E23").Value = '  +6.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.93%  '
$ws.Range("E27").Value = '  -6.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.58'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.58%  '
$ws.Range("E31").Value = '  -2.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '593.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.97'
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("D35").Value = '3.730.55'
$ws.Range("E35").Value = '  -6.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.91'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.37'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -9.37%  '
$ws.Range("E39").Value = '  +0.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.83'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.38%  '
$ws.Range("E41").Value = '  -5.51%  '
$ws.Range("E42").Value = '  -6.00%  '
$ws.Range("E43").Value = '  -2.80%  '
$ws.Range("D44").Value = '0.0₃0668'
$ws.Range("E44").Value = '  -5.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.26'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.82%  '
$ws.Range("E46").Value = '  -3.90%  '
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("E48").Value = '  -1.58%  '
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("E50").Value = '  -0.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '127.41'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.63%  '
